$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.118.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.83%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.479.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +5.21%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'555.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +5.52%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'181.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.93%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.637"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +8.92%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.471.66"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.20%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.01%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.631"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.65%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +14.54%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'54.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.14%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000275"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +6.90%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'9.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.76%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.055.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +5.43%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.488.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +5.38%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.121"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.81%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'18.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +6.24%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'66.202.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.25%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'11.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +7.42%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.991"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +4.10%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'417.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +10.47%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +9.76%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'85.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +5.65%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'4.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.75%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +7.15%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.69%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'12.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +9.40%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -2.01%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'9.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +11.10%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'30.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +4.73%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'6.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.21%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'620.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.26%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'11.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +5.32%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +5.29%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'60.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.34%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.147"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +17.78%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "'0.997"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.38%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "'37.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.63%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0₃0790"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +5.07%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.381"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.06%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +5.60%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'3.115.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +7.91%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.26%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +8.62%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.66%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0414"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +4.38%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'3.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +4.44%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.82%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.132"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +6.02%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'139.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.92%  "
$ws.Range("E51").Style = "Normal"
